$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New grade values for rows 2-12, columns B (2) through W (23).
# This rationalizes the fuzzy rule table used for grade computation.
$data = @(
    @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0.11, 0.18, 0.23, 0.25, 0.27, 0.28, 0.28, 0.28, 0.26, 0.24, 0.2, 0.12, 0),
    @(0.13, 0.13, 0.13, 0.13, 0.13, 0.13, 0.13, 0.13, 0.13, 0.13, 0.19, 0.23, 0.26, 0.28, 0.29, 0.29, 0.28, 0.27, 0.25, 0.21, 0.13, 0.12),
    @(0.23, 0.23, 0.23, 0.23, 0.23, 0.23, 0.23, 0.23, 0.23, 0.23, 0.23, 0.26, 0.28, 0.29, 0.3, 0.31, 0.3, 0.29, 0.27, 0.23, 0.21, 0.2),
    @(0.31, 0.31, 0.31, 0.31, 0.31, 0.31, 0.31, 0.31, 0.31, 0.31, 0.31, 0.31, 0.31, 0.33, 0.34, 0.34, 0.34, 0.32, 0.31, 0.27, 0.25, 0.24),
    @(0.36, 0.36, 0.36, 0.36, 0.36, 0.36, 0.36, 0.36, 0.36, 0.36, 0.36, 0.36, 0.36, 0.36, 0.37, 0.37, 0.37, 0.36, 0.32, 0.29, 0.27, 0.26),
    @(0.4, 0.4, 0.4, 0.4, 0.4, 0.4, 0.4, 0.4, 0.4, 0.4, 0.4, 0.4, 0.4, 0.4, 0.4, 0.4, 0.4, 0.37, 0.34, 0.3, 0.28, 0.28),
    @(0.41, 0.41, 0.41, 0.41, 0.41, 0.41, 0.41, 0.41, 0.41, 0.41, 0.41, 0.41, 0.41, 0.41, 0.41, 0.41, 0.4, 0.37, 0.34, 0.31, 0.29, 0.28),
    @(0.51, 0.51, 0.51, 0.51, 0.51, 0.51, 0.51, 0.51, 0.51, 0.51, 0.51, 0.51, 0.51, 0.51, 0.47, 0.41, 0.4, 0.37, 0.33, 0.3, 0.28, 0.27),
    @(0.61, 0.61, 0.61, 0.61, 0.61, 0.61, 0.61, 0.61, 0.61, 0.61, 0.61, 0.61, 0.58, 0.52, 0.47, 0.41, 0.4, 0.36, 0.31, 0.27, 0.25, 0.24),
    @(0.76, 0.76, 0.76, 0.76, 0.76, 0.76, 0.76, 0.76, 0.76, 0.76, 0.73, 0.65, 0.58, 0.52, 0.47, 0.41, 0.4, 0.36, 0.31, 0.23, 0.17, 0.16),
    @(1, 1, 1, 1, 1, 1, 1, 1, 1, 0.84, 0.73, 0.65, 0.58, 0.52, 0.47, 0.41, 0.4, 0.36, 0.31, 0.23, 0.13, 0)
)

$startRow = 2
$startCol = 2  # column B

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($startRow + $i, $startCol + $j).Value = $rowVals[$j]
    }
}
